$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 82  # was 81
$ws.Range("F5").Value = 348  # was 347
$ws.Range("F6").Value = 543  # was 542
$ws.Range("F8").Value = 11524  # was 11518
$ws.Range("F12").Value = 2085  # was 2082
$ws.Range("F16").Value = 227  # was 226
$ws.Range("F18").Value = 1176  # was 1175
$ws.Range("F19").Value = 151  # was 150
$ws.Range("F20").Value = 237  # was 236
$ws.Range("F21").Value = 724  # was 723
$ws.Range("F22").Value = 250  # was 132
$ws.Range("F23").Value = 252  # was 250
$ws.Range("F24").Value = 2395  # was 2394
$ws.Range("F26").Value = 3407  # was 3406
$ws.Range("F27").Value = 1047  # was 1041
$ws.Range("F29").Value = 10  # was 9
$ws.Range("F31").Value = 19  # was 18
$ws.Range("F32").Value = 963  # was 962
$ws.Range("F33").Value = 29  # was 28
$ws.Range("F36").Value = 11  # was 9
$ws.Range("F38").Value = 1820  # was 1718
$ws.Range("F39").Value = 4346  # was 4342
$ws.Range("F40").Value = 5435  # was 5433
$ws.Range("F43").Value = 26  # was 25
$ws.Range("F44").Value = 147  # was 145
$ws.Range("F45").Value = 250  # was 248
$ws.Range("F48").Value = 4087  # was 4088

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 4138  # was 4134
$ws.Range("F11").Value = 554  # was 545
$ws.Range("F14").Value = 3  # was 2

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 55  # was 54

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 55  # was 54
$ws.Range("F5").Value = 82  # was 81
$ws.Range("F7").Value = 348  # was 347
$ws.Range("F8").Value = 543  # was 542
$ws.Range("F9").Value = 11524  # was 11518
$ws.Range("F12").Value = 2085  # was 2082
$ws.Range("F15").Value = 227  # was 226
$ws.Range("F17").Value = 1176  # was 1175
$ws.Range("F18").Value = 151  # was 150
$ws.Range("F19").Value = 237  # was 236
$ws.Range("F20").Value = 4138  # was 4134
$ws.Range("F22").Value = 724  # was 723
$ws.Range("F23").Value = 251  # was 132
$ws.Range("F24").Value = 252  # was 250
$ws.Range("F26").Value = 1047  # was 1041
$ws.Range("F30").Value = 10  # was 9
$ws.Range("F32").Value = 19  # was 18
$ws.Range("F34").Value = 11  # was 9
$ws.Range("F38").Value = 147  # was 145
$ws.Range("F39").Value = 250  # was 248
$ws.Range("F43").Value = 4087  # was 4088
